$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.124.05'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '2.430.58'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '89.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.539'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.42%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.499'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0835'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '32.08'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.97%  '
$ws.Range('E12').Value = '  -1.87%  '
$ws.Range('D13').Value = '2.802.57'
$ws.Range('E13').Value = '  -1.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.73'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('D16').Value = '2.428.42'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.774'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').Value = '41.064.62'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.27'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.57'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.66'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.65%  '
$ws.Range('E24').Value = '  -1.90%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('E28').Value = '  -2.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.90'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.28'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0750'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.61%  '
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.99'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.91%  '
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.80'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.88%  '
$ws.Range('E40').Value = '  -2.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.93'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('D42').Value = '1.995.93'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('E45').Value = '  -3.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.10%  '
$ws.Range('D48').Value = '2.658.95'
$ws.Range('E48').Value = '  -1.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '95.28'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.85'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.31'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.59%  '
